# Refresh "想去人数" (interest/want-to-go counts) for the bilibili-scraped
# event rows that reappear on gh-pages re-generation (commit 456a3b4).
$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition listing)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F11").Value = 2065
$ws1.Range("F12").Value = 2065
$ws1.Range("F16").Value = 227
$ws1.Range("F18").Value = 4807
$ws1.Range("F19").Value = 132
$ws1.Range("F20").Value = 55
$ws1.Range("F34").Value = 794

# Sheet "演出" (performance listing)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F21").Value = 249

# Sheet "全部类型" (all-types aggregate listing)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F21").Value = 2065
$ws4.Range("F26").Value = 227
$ws4.Range("F28").Value = 4807
$ws4.Range("F29").Value = 55
$ws4.Range("F42").Value = 249
$ws4.Range("F44").Value = 794
